$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.632
$ws.Range("D5").Value = -8.063000000000001
$ws.Range("B8").Value = 5.369999999999999
$ws.Range("D8").Value = -7.85
$ws.Range("B10").Value = 6.052
$ws.Range("C11").Value = -12.312
$ws.Range("B12").Value = 5.356999999999999
$ws.Range("C12").Value = -11.659
$ws.Range("D12").Value = -7.629
$ws.Range("D13").Value = -7.677000000000001
$ws.Range("C15").Value = -13.685
$ws.Range("D15").Value = -8.078999999999999
$ws.Range("C17").Value = -13.271
$ws.Range("B18").Value = 5.241
$ws.Range("D21").Value = -7.713000000000001
$ws.Range("B25").Value = 5.340999999999999
$ws.Range("D25").Value = -7.51
$ws.Range("C26").Value = -12.545
$ws.Range("C27").Value = -12.661
$ws.Range("C28").Value = -12.418
$ws.Range("C32").Value = -11.726
$ws.Range("D32").Value = -7.331999999999999
$ws.Range("D36").Value = -7.327
$ws.Range("B37").Value = 8.398
$ws.Range("C37").Value = -11.862
$ws.Range("D38").Value = -7.795
$ws.Range("C41").Value = -12.06
$ws.Range("D41").Value = -8.164000000000001
$ws.Range("C47").Value = -12.608
$ws.Range("D50").Value = -7.975
$ws.Range("C51").Value = -12.916
$ws.Range("D52").Value = -8.077999999999999
$ws.Range("B55").Value = 4.864999999999999
$ws.Range("D59").Value = -7.722000000000001
$ws.Range("C65").Value = -12.161
$ws.Range("D67").Value = -7.767
$ws.Range("B68").Value = 4.755
$ws.Range("C73").Value = -12.347
$ws.Range("B77").Value = 6.433
$ws.Range("B78").Value = 8.016999999999999
$ws.Range("B79").Value = 5.846
$ws.Range("B80").Value = 8.204000000000001
$ws.Range("B81").Value = 6.187
$ws.Range("B82").Value = 5.579
$ws.Range("B84").Value = 5.584000000000001
$ws.Range("C84").Value = -12.9
$ws.Range("D84").Value = -7.883000000000001
$ws.Range("C85").Value = -12.573
$ws.Range("D86").Value = -7.874000000000001
$ws.Range("D88").Value = -7.937
$ws.Range("C89").Value = -13.235
$ws.Range("D89").Value = -8.159000000000001
$ws.Range("C93").Value = -12.586
$ws.Range("C95").Value = -12.2
$ws.Range("D95").Value = -7.736
$ws.Range("C98").Value = -12.959
$ws.Range("C99").Value = -11.784
$ws.Range("B101").Value = 6.007
$ws.Range("C101").Value = -12.822
$ws.Range("B102").Value = 6.449
$ws.Range("C102").Value = -12.919
$ws.Range("D105").Value = -7.681

Write-Host "Applied 62 cell updates"
